$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.182.92"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.586.89"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.07%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "2.597.27"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("D14").Value = "3.049.64"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "59.094.20"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "2.598.50"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.404"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "0.0₃0722"
$ws.Range("E30").Value = "  -3.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.78%  "
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.812"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.57%  "
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.602"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "271.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.45%  "
$ws.Range("D49").Value = "1.969.43"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.83%  "
